$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "51.978.61"
$ws.Range("E2").Value = "  -0.37%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.964.12"
$ws.Range("E3").Value = "  +3.01%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - BNB
$ws.Range("D5").Value = "353.59"
$ws.Range("E5").Value = "  +0.25%  "

# Row 6 - Solana
$ws.Range("D6").Value = "112.30"
$ws.Range("E6").Value = "  -0.53%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +1.05%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.03%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +1.62%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "39.74"
$ws.Range("E10").Value = "  -2.12%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "0.0899"
$ws.Range("E11").Value = "  +5.22%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +0.98%  "

# Row 13 - Chainlink
$ws.Range("D13").Value = "19.99"
$ws.Range("E13").Value = "  -1.24%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "8.02"
$ws.Range("E14").Value = "  +2.22%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.437.03"
$ws.Range("E15").Value = "  +3.17%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.968.81"
$ws.Range("E16").Value = "  +1.82%  "

# Row 17 - Polygon
$ws.Range("D17").Value = "1.00"
$ws.Range("E17").Value = "  +0.89%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "52.093.79"
$ws.Range("E18").Value = "  -0.13%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "7.74"
$ws.Range("E19").Value = "  +0.61%  "

# Row 20 - InternetComputer(DFINITY)
$ws.Range("D20").Value = "14.55"
$ws.Range("E20").Value = "  +6.20%  "

# Row 21 - ImmutableX
$ws.Range("D21").Value = "3.33"
$ws.Range("E21").Value = "  -2.15%  "

# Row 22 - ShibaInu
$ws.Range("D22").Value = "0.0₃0993"
$ws.Range("E22").Value = "  +1.63%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "71.51"
$ws.Range("E23").Value = "  +1.02%  "

# Row 24 - BitcoinCash
$ws.Range("D24").Value = "271.07"
$ws.Range("E24").Value = "  +0.05%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +0.62%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  +9.77%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "27.48"
$ws.Range("E27").Value = "  +3.45%  "

# Row 28 - Filecoin
$ws.Range("D28").Value = "7.61"
$ws.Range("E28").Value = "  +20.58%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  +0.01%  "

# Row 30 - Hedera
$ws.Range("D30").Value = "0.110"
$ws.Range("E30").Value = "  +22.54%  "

# Row 31 - Cosmos
$ws.Range("D31").Value = "10.76"
$ws.Range("E31").Value = "  +1.77%  "

# Row 32 - InjectiveProtocol
$ws.Range("D32").Value = "37.86"
$ws.Range("E32").Value = "  -2.74%  "

# Row 33 - RenderToken
$ws.Range("D33").Value = "6.21"
$ws.Range("E33").Value = "  +10.13%  "

# Row 34 - OKB
$ws.Range("D34").Value = "53.23"
$ws.Range("E34").Value = "  +1.08%  "

# Row 35 & 36 - VeChain and Toncoin swap places
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").Value = "0.0452"
$ws.Range("E35").Value = "  -0.78%  "

$ws.Range("B36").Value = "Toncoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D36").Value = "1.98"
$ws.Range("E36").Value = "  -12.95%  "

# Row 37 - FirstDigitalUSD
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.15%  "

# Row 38 - LidoDAOToken
$ws.Range("E38").Value = "  +3.55%  "

# Row 39 - Celestia
$ws.Range("D39").Value = "19.02"
$ws.Range("E39").Value = "  +0.17%  "

# Row 40 - ARBITRUM
$ws.Range("D40").Value = "2.07"
$ws.Range("E40").Value = "  +1.47%  "

# Row 41 - Stacks
$ws.Range("D41").Value = "2.70"
$ws.Range("E41").Value = "  +4.08%  "

# Row 42 - EnergySwap
$ws.Range("E42").Value = "  +6.06%  "

# Row 43 - Stellar
$ws.Range("E43").Value = "  +1.65%  "

# Row 44 - WEMIXToken
$ws.Range("E44").Value = "  -2.28%  "

# Row 45 - NEARProtocol
$ws.Range("D45").Value = "3.59"
$ws.Range("E45").Value = "  +0.78%  "

# Row 46 - ApeXProtocol
$ws.Range("E46").Value = "  +1.78%  "

# Row 47 - Maker
$ws.Range("D47").Value = "2.183.98"
$ws.Range("E47").Value = "  +0.06%  "

# Row 48 - Monero
$ws.Range("D48").Value = "113.77"
$ws.Range("E48").Value = "  -7.22%  "

# Row 49 - TheGraph
$ws.Range("E49").Value = "  -0.44%  "

# Row 50 - BEAM
$ws.Range("D50").Value = "0.0340"
$ws.Range("E50").Value = "  +5.58%  "

# Row 51 - SEI
$ws.Range("D51").Value = "0.940"
$ws.Range("E51").Value = "  -2.43%  "
